# Auto-generated edit script: update 2022 (column I) violent crime counts
# per commit "Add data for 2022-12-04"

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 6783
$ws.Range("I3").Value = 7062
$ws.Range("I4").Value = 1624
$ws.Range("I6").Value = 8208
$ws.Range("I7").Value = 24333

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 189
$ws.Range("I5").Value = 74
$ws.Range("I7").Value = 766
$ws.Range("I8").Value = 1453
$ws.Range("I10").Value = 179
$ws.Range("I15").Value = 283
$ws.Range("I18").Value = 191
$ws.Range("I20").Value = 603
$ws.Range("I21").Value = 107
$ws.Range("I22").Value = 69
$ws.Range("I26").Value = 33
$ws.Range("I28").Value = 13
$ws.Range("I29").Value = 1464
$ws.Range("I31").Value = 246
$ws.Range("I33").Value = 1077
$ws.Range("I39").Value = 18
$ws.Range("I41").Value = 106
$ws.Range("I42").Value = 893
$ws.Range("I44").Value = 185
$ws.Range("I47").Value = 176
$ws.Range("I52").Value = 555
$ws.Range("I53").Value = 269
$ws.Range("I54").Value = 484
$ws.Range("I55").Value = 281
$ws.Range("I58").Value = 15
$ws.Range("I59").Value = 41
$ws.Range("I60").Value = 138
$ws.Range("I63").Value = 78
$ws.Range("I64").Value = 193
$ws.Range("I65").Value = 568
$ws.Range("I66").Value = 69
$ws.Range("I67").Value = 926
$ws.Range("I70").Value = 42
$ws.Range("I76").Value = 348
$ws.Range("I79").Value = 697
$ws.Range("I80").Value = 77
$ws.Range("I83").Value = 525
$ws.Range("I85").Value = 1090
$ws.Range("I86").Value = 154
$ws.Range("I90").Value = 314
$ws.Range("I93").Value = 138
$ws.Range("I94").Value = 248
$ws.Range("I95").Value = 370
$ws.Range("I96").Value = 281
$ws.Range("I97").Value = 211
$ws.Range("I98").Value = 170
$ws.Range("I101").Value = 24333

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 311
$ws.Range("I4").Value = 49
$ws.Range("I6").Value = 285
$ws.Range("I7").Value = 1090

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I6").Value = 181
$ws.Range("I7").Value = 555

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 430
$ws.Range("I3").Value = 421
$ws.Range("I6").Value = 466
$ws.Range("I7").Value = 1453

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I4").Value = 23
$ws.Range("I6").Value = 128
$ws.Range("I7").Value = 269

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I2").Value = 249
$ws.Range("I3").Value = 234
$ws.Range("I6").Value = 207
$ws.Range("I7").Value = 766

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I2").Value = 82
$ws.Range("I7").Value = 281

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I2").Value = 219
$ws.Range("I3").Value = 346
$ws.Range("I7").Value = 926

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("I3").Value = 62
$ws.Range("I6").Value = 95
$ws.Range("I7").Value = 246

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I4").Value = 23
$ws.Range("I6").Value = 174
$ws.Range("I7").Value = 568

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I3").Value = 192
$ws.Range("I7").Value = 525

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I2").Value = 131
$ws.Range("I3").Value = 131
$ws.Range("I7").Value = 370

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 243
$ws.Range("I7").Value = 1077

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I2").Value = 103
$ws.Range("I7").Value = 484

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 432
$ws.Range("I6").Value = 403
$ws.Range("I7").Value = 1464

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("I2").Value = 60
$ws.Range("I6").Value = 53
$ws.Range("I7").Value = 185

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I2").Value = 71
$ws.Range("I4").Value = 37
$ws.Range("I7").Value = 348

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("I2").Value = 32
$ws.Range("I7").Value = 106

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I6").Value = 335
$ws.Range("I7").Value = 893

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("I3").Value = 36
$ws.Range("I7").Value = 179

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I3").Value = 88
$ws.Range("I6").Value = 86
$ws.Range("I7").Value = 281

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("I6").Value = 80
$ws.Range("I7").Value = 107

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I2").Value = 203
$ws.Range("I3").Value = 227
$ws.Range("I6").Value = 200
$ws.Range("I7").Value = 697

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("I6").Value = 63
$ws.Range("I7").Value = 193

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I2").Value = 169
$ws.Range("I6").Value = 208
$ws.Range("I7").Value = 603

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("I3").Value = 42
$ws.Range("I4").Value = 6
$ws.Range("I7").Value = 191

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("I2").Value = 39
$ws.Range("I7").Value = 138

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I2").Value = 48
$ws.Range("I6").Value = 142
$ws.Range("I7").Value = 248

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I6").Value = 58
$ws.Range("I7").Value = 176

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I2").Value = 82
$ws.Range("I3").Value = 66
$ws.Range("I6").Value = 108
$ws.Range("I7").Value = 283

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I6").Value = 108
$ws.Range("I7").Value = 170

$ws = $wb.Worksheets.Item('East Village')
$ws.Range("I6").Value = 21
$ws.Range("I7").Value = 33

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range("I5").Value = 10
$ws.Range("I6").Value = 18

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("I6").Value = 30
$ws.Range("I7").Value = 69

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range("I3").Value = 5
$ws.Range("I7").Value = 41

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I2").Value = 70
$ws.Range("I7").Value = 189

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("I2").Value = 34
$ws.Range("I6").Value = 137
$ws.Range("I7").Value = 211

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("I2").Value = 13
$ws.Range("I7").Value = 42

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("I4").Value = 4
$ws.Range("I7").Value = 74

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("I4").Value = 74
$ws.Range("I7").Value = 154

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I2").Value = 102
$ws.Range("I7").Value = 314

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("I2").Value = 48
$ws.Range("I7").Value = 138

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("I6").Value = 19
$ws.Range("I7").Value = 69

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("I6").Value = 43
$ws.Range("I7").Value = 77

$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range("I4").Value = 2
$ws.Range("I7").Value = 15

$ws = $wb.Worksheets.Item('Edison Park')
$ws.Range("I3").Value = 5
$ws.Range("I7").Value = 13

Write-Output "Updated 165 cells across 47 sheets"
